# Update cryptocurrency price/volume data to reflect the Dec 1 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$text)
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = $origStyle
}

$ws.Range('D2').Value = '38.801.76'
$ws.Range('E2').Value = '  +2.72%  '

$ws.Range('D3').Value = '2.093.38'
$ws.Range('E3').Value = '  +2.65%  '

$ws.Range('E4').Value = '  -0.05%  '

Set-TextValue $ws.Range('D5') '228.12'
$ws.Range('E5').Value = '  +0.25%  '

Set-TextValue $ws.Range('D6') '0.612'
$ws.Range('E6').Value = '  +0.90%  '

Set-TextValue $ws.Range('D7') '60.74'
$ws.Range('E7').Value = '  +2.12%  '

Set-TextValue $ws.Range('D8') '0.999'
$ws.Range('E8').Value = '  -0.08%  '

Set-TextValue $ws.Range('D10') '0.0834'
$ws.Range('E10').Value = '  -0.86%  '

$ws.Range('E11').Value = '  -0.28%  '

$ws.Range('D12').Value = '2.401.04'
$ws.Range('E12').Value = '  +2.54%  '

Set-TextValue $ws.Range('D13') '14.95'
$ws.Range('E13').Value = '  +3.43%  '

Set-TextValue $ws.Range('D14') '22.02'

Set-TextValue $ws.Range('D15') '0.797'
$ws.Range('E15').Value = '  +3.10%  '

$ws.Range('E16').Value = '  -0.28%  '

$ws.Range('D17').Value = '2.088.71'
$ws.Range('E17').Value = '  +2.41%  '

$ws.Range('D18').Value = '38.673.49'
$ws.Range('E18').Value = '  +2.40%  '

Set-TextValue $ws.Range('D19') '72.02'
$ws.Range('E19').Value = '  +3.62%  '

$ws.Range('E20').Value = '  +2.32%  '

$ws.Range('D21').Value = '0.0₃0834'
$ws.Range('E21').Value = '  +1.31%  '

Set-TextValue $ws.Range('D22') '226.33'
$ws.Range('E22').Value = '  +1.22%  '

$ws.Range('E23').Value = '  -0.48%  '

Set-TextValue $ws.Range('D24') '2.48'
$ws.Range('E24').Value = '  +1.99%  '

Set-TextValue $ws.Range('D25') '2.34'
$ws.Range('E25').Value = '  +2.78%  '

Set-TextValue $ws.Range('D26') '170.63'
$ws.Range('E26').Value = '  +1.37%  '

Set-TextValue $ws.Range('D27') '9.48'
$ws.Range('E27').Value = '  +1.09%  '

Set-TextValue $ws.Range('D28') '0.137'
$ws.Range('E28').Value = '  +6.44%  '

Set-TextValue $ws.Range('D29') '1.41'
$ws.Range('E29').Value = '  +10.40%  '

Set-TextValue $ws.Range('D30') '19.16'
$ws.Range('E30').Value = '  +1.97%  '

$ws.Range('E31').Value = '  +0.44%  '

$ws.Range('E32').Value = '  +4.53%  '

$ws.Range('E33').Value = '  +5.40%  '

$ws.Range('E34').Value = '  +2.13%  '

Set-TextValue $ws.Range('D35') '0.0612'
$ws.Range('E35').Value = '  +1.22%  '

Set-TextValue $ws.Range('D36') '6.42'
$ws.Range('E36').Value = '  +0.00%  '

$ws.Range('E37').Value = '  +2.04%  '

$ws.Range('E38').Value = '  +2.83%  '

$ws.Range('E39').Value = '  -0.04%  '

Set-TextValue $ws.Range('D40') '18.31'
$ws.Range('E40').Value = '  +1.42%  '

$ws.Range('D41').Value = '1.536.64'
$ws.Range('E41').Value = '  +0.61%  '

Set-TextValue $ws.Range('D42') '100.91'
$ws.Range('E42').Value = '  +3.72%  '

$ws.Range('E43').Value = '  +3.03%  '

$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D44') '0.0925'
$ws.Range('E44').Value = '  +2.32%  '

$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D45') '2.82'
$ws.Range('E45').Value = '  -0.89%  '

Set-TextValue $ws.Range('D46') '7.65'
$ws.Range('E46').Value = '  +9.36%  '

$ws.Range('E47').Value = '  +0.62%  '

Set-TextValue $ws.Range('D48') '4.10'
$ws.Range('E48').Value = '  -3.06%  '

$ws.Range('E49').Value = '  +2.94%  '

$ws.Range('E50').Value = '  +1.16%  '

$ws.Range('E51').Value = '  +2.59%  '

